$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.884.29"
$ws.Range("D3").Value = "3.277.38"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.68%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("E9").Value = "  +7.59%  "
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("E11").Value = "  +5.94%  "
$ws.Range("D12").Value = "3.850.50"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.138"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.39%  "
$ws.Range("D15").Value = "67.843.79"
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").Value = "3.286.05"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("E19").Value = "  +4.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "377.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.48%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.72%  "
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000120"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.32%  "
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("E27").Value = "  +3.09%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  +5.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.96%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.05%  "
$ws.Range("E34").Value = "  +5.24%  "
$ws.Range("E35").Value = "  +5.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("E40").Value = "  +10.60%  "
$ws.Range("E41").Value = "  +10.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "351.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "2.656.88"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0284"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("E50").Value = "  +5.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.12"
$ws.Range("D51").Style = "Normal"
